$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sets a cell value as literal text (preventing Excel's automatic
# numeric/percentage coercion), matching the inline-string cell type
# already used throughout this sheet, and resets the style index so
# no incidental number-format/quote-prefix style is introduced.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "305.30"
Set-TextValue "E2" "0.19%"
Set-TextValue "D3" "35.97"
Set-TextValue "E3" "0.56%"
Set-TextValue "D4" "5.055"
Set-TextValue "E4" "-0.26%"
Set-TextValue "D5" "0.08046"
Set-TextValue "E5" "-0.04%"
Set-TextValue "D6" "1.869"
Set-TextValue "E6" "-2.68%"
Set-TextValue "D7" "4.120"
Set-TextValue "E7" "-1.07%"
Set-TextValue "D8" "7.785"
Set-TextValue "E8" "-0.74%"
Set-TextValue "D9" "0.9263"
Set-TextValue "E9" "-0.48%"
Set-TextValue "D10" "0.1363"
Set-TextValue "E10" "1.64%"
Set-TextValue "D11" "0.1896"
Set-TextValue "E11" "-0.05%"
Set-TextValue "D12" "0.09061"
Set-TextValue "E12" "-0.97%"
Set-TextValue "D13" "0.03431"
Set-TextValue "E13" "-1.26%"
Set-TextValue "D14" "0.09891"
Set-TextValue "E14" "-0.15%"
Set-TextValue "D15" "0.001404"
Set-TextValue "E15" "-0.97%"
Set-TextValue "D16" "0.006090"
Set-TextValue "E16" "-8.72%"
Set-TextValue "E17" "6.31%"
Set-TextValue "D18" "3.390"
Set-TextValue "E18" "13.42%"
Set-TextValue "D19" "0.3418"
Set-TextValue "E19" "-0.16%"
Set-TextValue "D20" "0.1333"
Set-TextValue "E20" "-0.30%"
Set-TextValue "D21" "4.814"
Set-TextValue "E21" "-7.20%"
Set-TextValue "D22" "0.2384"
Set-TextValue "E22" "-6.05%"
Set-TextValue "D23" "0.04356"
Set-TextValue "E23" "-1.42%"
Set-TextValue "E24" "-0.62%"
Set-TextValue "D25" "0.004287"
Set-TextValue "E25" "-8.84%"
Set-TextValue "E27" "-0.41%"
Set-TextValue "E28" "41.88%"
Set-TextValue "D39" "0.02001"
Set-TextValue "E39" "0.25%"
Set-TextValue "D40" "0.05105"
Set-TextValue "E40" "-0.30%"
Set-TextValue "D41" "0.007515"
Set-TextValue "D42" "0.01007"
Set-TextValue "E42" "-6.58%"
Set-TextValue "D43" "0.1357"
Set-TextValue "E43" "-0.48%"
Set-TextValue "D44" "0.002157"
Set-TextValue "E44" "2.43%"
Set-TextValue "D45" "0.009619"
Set-TextValue "E45" "-10.51%"
Set-TextValue "D46" "0.00006221"
Set-TextValue "E46" "-1.41%"
Set-TextValue "D47" "0.00000000749"
Set-TextValue "E47" "-0.32%"
Set-TextValue "E48" "-0.16%"
Set-TextValue "D49" "0.001249"
Set-TextValue "E49" "-22.12%"
Set-TextValue "D50" "0.00002098"
Set-TextValue "E50" "-0.32%"
Set-TextValue "D51" "0.0001998"
Set-TextValue "E51" "-0.32%"
